# Atualizacao de bases das ligas, do dia: 24-02-2024 as 23:13
# Swap the data (columns B:AC) between row pairs (189,191) and (207,208),
# keeping column A (the row index) unchanged for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Algeria Division 1")

$firstCol = 2
$lastCol = 29

$rangeA = $ws.Range($ws.Cells.Item(189, $firstCol), $ws.Cells.Item(189, $lastCol))
$rangeB = $ws.Range($ws.Cells.Item(191, $firstCol), $ws.Cells.Item(191, $lastCol))

$valuesA = $rangeA.Value2
$valuesB = $rangeB.Value2

$rangeA.Value2 = $valuesB
$rangeB.Value2 = $valuesA

$rangeC = $ws.Range($ws.Cells.Item(207, $firstCol), $ws.Cells.Item(207, $lastCol))
$rangeD = $ws.Range($ws.Cells.Item(208, $firstCol), $ws.Cells.Item(208, $lastCol))

$valuesC = $rangeC.Value2
$valuesD = $rangeD.Value2

$rangeC.Value2 = $valuesD
$rangeD.Value2 = $valuesC
